$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": add a new column Q ("30-jun") with header + 24 values
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the last existing header/data column (P) into the
# new column Q so the new cells keep the same style as their neighbours.
$wsSpot.Range("P1:P25").Copy($wsSpot.Range("Q1:Q25"))

$wsSpot.Range("Q1").Value = "30-jun"

$qValues = @{
    2  = 98.73
    3  = 89.03
    4  = 85.03
    5  = 84.13
    6  = 82.98
    7  = 82.02
    8  = 97.73
    9  = 111.68
    10 = 108.58
    11 = 92.42
    12 = 79.83
    13 = 69.05
    14 = 52.44
    15 = 46.31
    16 = 51.53
    17 = 64.76000000000001
    18 = 83.83
    19 = 93.56999999999999
    20 = 110.65
    21 = 178.94
    22 = 185
    23 = 175.01
    24 = 157
    25 = 122.65
}

foreach ($row in $qValues.Keys) {
    $wsSpot.Cells.Item($row, 17).Value = $qValues[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append two new rows (14, 15) with dates + last price
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date-like strings to be stored as plain text (matching the
# existing "Date" column) instead of being auto-converted into date serial
# numbers by Excel, then clear the formatting change so the cell keeps the
# default (unstyled) look of the other date cells in the column.
$wsGaz.Range("A14").NumberFormat = "@"
$wsGaz.Range("A14").Value = "2025-06-28"
$wsGaz.Range("A14").ClearFormats()
$wsGaz.Range("B14").Value = 32.675

$wsGaz.Range("A15").NumberFormat = "@"
$wsGaz.Range("A15").Value = "2025-06-29"
$wsGaz.Range("A15").ClearFormats()
$wsGaz.Range("B15").Value = 32.675

# ---------------------------------------------------------------------------
# Sheet "CO2": append two new rows (14, 15) with dates + last price
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A14").NumberFormat = "@"
$wsCo2.Range("A14").Value = "2025-06-28"
$wsCo2.Range("A14").ClearFormats()
$wsCo2.Range("B14").Value = 69.92

$wsCo2.Range("A15").NumberFormat = "@"
$wsCo2.Range("A15").Value = "2025-06-29"
$wsCo2.Range("A15").ClearFormats()
$wsCo2.Range("B15").Value = 69.92

Write-Host "Edit complete"
